$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 12131.429
$ws.Range("I18").Value = 12358.333
$ws.Range("K18").Value = 12358.333
$ws.Range("M18").Value = -12074.333

$ws.Range("H38").Value = 379.4737
$ws.Range("I38").Value = 158.57143
$ws.Range("K38").Value = 475.71429
$ws.Range("M38").Value = -103.71429

$ws.Range("H43").Value = 11999
$ws.Range("I43").Value = 11999
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 11999
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -11930
$ws.Range("N43").ClearContents()

$ws.Range("H70").Value = 7000.2144
$ws.Range("I70").Value = 7230.769
$ws.Range("J70").Value = 4003
$ws.Range("K70").Value = 21692.307
$ws.Range("L70").Value = 12009
$ws.Range("M70").Value = -21422.307
$ws.Range("N70").Value = -12549

$ws.Range("H73").Value = 7000.2144
$ws.Range("I73").Value = 7230.769
$ws.Range("J73").Value = 4003
$ws.Range("K73").Value = 21692.307
$ws.Range("L73").Value = 12009
$ws.Range("M73").Value = -20756.307
$ws.Range("N73").Value = -13881

$ws.Range("H92").Value = 696.4666999999999
$ws.Range("I92").Value = 704.4545000000001
$ws.Range("J92").Value = 674.5
$ws.Range("K92").Value = 704.4545000000001
$ws.Range("L92").Value = 674.5
$ws.Range("M92").Value = 543.5454999999999
$ws.Range("N92").Value = -3170.5

$ws.Range("H96").Value = 2028.5
$ws.Range("J96").Value = 1675
$ws.Range("L96").Value = 5025
$ws.Range("N96").Value = -7771

$ws.Range("H98").Value = 691.2727
$ws.Range("I98").Value = 691.2727
$ws.Range("K98").Value = 691.2727
$ws.Range("M98").Value = 806.7273

$ws.Range("H107").Value = 2227
$ws.Range("I107").Value = 2280
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 2280
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = -360
$ws.Range("N107").Value = -5590

$ws.Range("H113").Value = 2524.75
$ws.Range("I113").Value = 2099.8333
$ws.Range("J113").Value = 3799.5
$ws.Range("K113").Value = 2099.8333
$ws.Range("L113").Value = 3799.5
$ws.Range("M113").Value = 1154.1667
$ws.Range("N113").Value = -10307.5

$ws.Range("H122").Value = 691.2727
$ws.Range("I122").Value = 691.2727
$ws.Range("K122").Value = 2073.8181
$ws.Range("M122").Value = 376.1819

$ws.Range("H130").Value = 64999
$ws.Range("J130").Value = 64999
$ws.Range("L130").Value = 64999
$ws.Range("N130").Value = -75039

$ws.Range("H137").Value = 2280.8
$ws.Range("I137").Value = 2101
$ws.Range("K137").Value = 6303
$ws.Range("M137").Value = -3753

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5136.1665
$ws.Range("I28").Value = 5136.1665
$ws.Range("K28").Value = 5136.1665
$ws.Range("M28").Value = -4944.1665

$ws.Range("H32").Value = 374146.47
$ws.Range("I32").Value = 1652.8
$ws.Range("K32").Value = 1652.8
$ws.Range("M32").Value = -1365.8

$ws.Range("H45").Value = 1747.7059
$ws.Range("I45").Value = 1433.25
$ws.Range("K45").Value = 1433.25
$ws.Range("M45").Value = -1056.25

$ws.Range("H99").Value = 5136.1665
$ws.Range("I99").Value = 5136.1665
$ws.Range("K99").Value = 5136.1665
$ws.Range("M99").Value = -2141.1665

$ws.Range("H122").Value = 2778.8823
$ws.Range("I122").Value = 1742.7
$ws.Range("K122").Value = 5228.1
$ws.Range("M122").Value = -2778.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 400
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -746

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H99").Value = 1230.2858
$ws.Range("I99").Value = 1206.8
$ws.Range("K99").Value = 1206.8
$ws.Range("M99").Value = 291.2

$ws.Range("H107").Value = 2156.1428
$ws.Range("I107").Value = 1926
$ws.Range("K107").Value = 1926
$ws.Range("M107").Value = -6

$ws.Range("H134").Value = 2598.9092
$ws.Range("I134").Value = 2598.9092
$ws.Range("K134").Value = 7796.7276
$ws.Range("M134").Value = -5261.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2383.3333
$ws.Range("I7").Value = 2306
$ws.Range("J7").Value = 3002
$ws.Range("K7").Value = 2306
$ws.Range("L7").Value = 3002
$ws.Range("M7").Value = -2193
$ws.Range("N7").Value = -3228

$ws.Range("H31").Value = 5839.0586
$ws.Range("I31").Value = 2926.2856
$ws.Range("K31").Value = 2926.2856
$ws.Range("M31").Value = -2631.2856

$ws.Range("H34").Value = 5839.0586
$ws.Range("I34").Value = 2926.2856
$ws.Range("K34").Value = 2926.2856
$ws.Range("M34").Value = -2724.2856

$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246

$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232

$ws.Range("H122").Value = 2727.7
$ws.Range("I122").Value = 1756
$ws.Range("K122").Value = 5268
$ws.Range("M122").Value = -2818

$ws.Range("H132").Value = 5825.1665
$ws.Range("I132").Value = 5077.385
$ws.Range("K132").Value = 15232.155
$ws.Range("M132").Value = -12702.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H51").Value = 333.25
$ws.Range("I51").Value = 100
$ws.Range("J51").Value = 566.5
$ws.Range("K51").Value = 300
$ws.Range("L51").Value = 1699.5
$ws.Range("M51").Value = 160
$ws.Range("N51").Value = -2619.5

$ws.Range("H57").Value = 8499
$ws.Range("I57").Value = 8249
$ws.Range("K57").Value = 24747
$ws.Range("M57").Value = -24188

$ws.Range("H113").Value = 788
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 860
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2580
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6920

$ws.Range("H137").Value = 4458.5
$ws.Range("I137").Value = 5000
$ws.Range("K137").Value = 15000
$ws.Range("M137").Value = -9900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 34746.57
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 34746.57
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 34746.57
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -35236.57

$ws.Range("H24").Value = 111618.13
$ws.Range("I24").Value = 152384.28
$ws.Range("K24").Value = 152384.28
$ws.Range("M24").Value = -152211.28

$ws.Range("H132").Value = 5846.8335
$ws.Range("I132").Value = 4355.6665
$ws.Range("J132").Value = 7338
$ws.Range("K132").Value = 13066.9995
$ws.Range("L132").Value = 22014
$ws.Range("M132").Value = -10536.9995
$ws.Range("N132").Value = -27074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1331.5
$ws.Range("I22").Value = 1094.5
$ws.Range("J22").Value = 1450
$ws.Range("K22").Value = 1094.5
$ws.Range("L22").Value = 1450
$ws.Range("M22").Value = -799.5
$ws.Range("N22").Value = -2040

$ws.Range("H27").Value = 1331.5
$ws.Range("I27").Value = 1094.5
$ws.Range("J27").Value = 1450
$ws.Range("K27").Value = 1094.5
$ws.Range("L27").Value = 1450
$ws.Range("M27").Value = -987.5
$ws.Range("N27").Value = -1664

$ws.Range("H42").Value = 5000000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H46").Value = 2250
$ws.Range("J46").Value = 2225
$ws.Range("L46").Value = 2225
$ws.Range("N46").Value = -2601

$ws.Range("H49").Value = 5000000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H122").Value = 3537.6897
$ws.Range("I122").Value = 3012.375
$ws.Range("J122").Value = 3737.8096
$ws.Range("K122").Value = 9037.125
$ws.Range("L122").Value = 11213.4288
$ws.Range("M122").Value = -6587.125
$ws.Range("N122").Value = -16113.4288

$ws.Range("H132").Value = 262034.25
$ws.Range("I132").Value = 336045.66
$ws.Range("K132").Value = 1008136.98
$ws.Range("M132").Value = -1005606.98

$ws.Range("H136").Value = 5774.5
$ws.Range("I136").Value = 5999.3335
$ws.Range("K136").Value = 17998.0005
$ws.Range("M136").Value = -15448.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 21979.6
$ws.Range("I9").Value = 19974.5
$ws.Range("J9").Value = 30000
$ws.Range("K9").Value = 19974.5
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = -19834.5
$ws.Range("N9").Value = -30280

$ws.Range("H126").Value = 1744.5
$ws.Range("I126").Value = 1918.6
$ws.Range("K126").Value = 5755.799999999999
$ws.Range("M126").Value = -3285.799999999999

$ws.Range("H130").Value = 24500
$ws.Range("J130").Value = 24500
$ws.Range("L130").Value = 24500
$ws.Range("N130").Value = -34540

$ws.Range("H132").Value = 2974.889
$ws.Range("I132").Value = 2695.8572
$ws.Range("K132").Value = 8087.571599999999
$ws.Range("M132").Value = -5557.571599999999
